$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 63: Recursion / Scaler / Implement pow(...) / Solution60 ---
$ws.Range("A63").Value = "Recursion"
$ws.Range("B63").Value = "Scaler"
$ws.Range("C63").WrapText = $true
$ws.Range("C63").Value = "Implement pow(A, B) % C.`nIn other words, given A, B and C, Find (AB % C)."
$ws.Range("D63").Value = "Solution60"
$ws.Rows.Item(63).RowHeight = 30

# --- Row 64: Sorting / Scaler / Insertion Sort / Solution62 ---
$ws.Range("A64").Value = "Sorting"
$ws.Range("B64").Value = "Scaler"
$ws.Range("C64").WrapText = $true
$ws.Range("C64").Value = "Insertion Sort"
$ws.Range("D64").Value = "Solution62"

# --- Row 65: Array / Others / Find max product.../ Solution61 ---
$ws.Range("A65").Value = "Array"
$ws.Range("B65").Value = "Others"
$ws.Range("D65").Value = "Solution61"
$ws.Range("C65").WrapText = $true
$ws.Range("C65").Value = "Find max product of sub Array from a given Array"

# --- Row 66: Sorting / Scaler / Selection Sort / Solution63 ---
$ws.Range("A66").Value = "Sorting"
$ws.Range("B66").Value = "Scaler"
$ws.Range("C66").WrapText = $true
$ws.Range("C66").Value = "Selection Sort"
$ws.Range("D66").Value = "Solution63"

# --- Remove the AutoFilter that used to sit on A2:A57 ---
$ws.AutoFilterMode = $false

# --- Extend the Topic data-validation list to include "Sorting" ---
$topicValidation = $ws.Range("A1:A1048576").Validation
$topicValidation.Formula1 = '"Prime number,Subsequence,Tree,Stack,Queue,Recursion,Hashing,Array, String,Bit Manupulation,Loop,Maths,Modulus,Sorting"'

# --- Move the selection to where the author ended up (C62) ---
$null = $ws.Range("C62").Select()
